# Track-it Kravspecifikationer: fill in the "Implementeret / Testet"
# columns for the newly-handled SVIPT/ADAM requirement rows.
#
# Table 4 ("Feature | Implementeret | Testet") rows (1-based, row 1 is the
# header row):
#   18 - "Gør cursor 1-2 px større"
#   19 - "I SVIPT ændre cursor farve til det target / baseline ..."
#   20 - "Undersøg hvorfor programmet ikke lukker ordentligt ned ..."
#   21 - "Optimer forbindelse til ADAM"

$d = $word.ActiveDocument
$t = $d.Tables.Item(4)

$t.Cell(18, 2).Range.Text = "X"
$t.Cell(18, 3).Range.Text = "X"

$t.Cell(19, 2).Range.Text = "X"
$t.Cell(19, 3).Range.Text = "X"

$t.Cell(20, 2).Range.Text = "Skal undersøges på Lab3 pc"

$t.Cell(21, 2).Range.Text = "X"
$t.Cell(21, 3).Range.Text = "X"
